$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.516.37'
$ws.Range('E2').Value = '  +0.80%  '

$ws.Range('D3').Value = '2.982.95'
$ws.Range('E3').Value = '  +2.47%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.31'
$ws.Range('E5').Value = '  +3.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.78'
$ws.Range('E6').Value = '  +2.50%  '

$ws.Range('E7').Value = '  +0.39%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +1.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.42'
$ws.Range('E10').Value = '  +1.94%  '

$ws.Range('E11').Value = '  +0.35%  '

$ws.Range('E12').Value = '  +1.19%  '

$ws.Range('D13').Value = '3.455.51'
$ws.Range('E13').Value = '  +2.33%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.47'
$ws.Range('E14').Value = '  +0.38%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.56'
$ws.Range('E15').Value = '  +2.51%  '

$ws.Range('D16').Value = '2.985.31'
$ws.Range('E16').Value = '  +2.34%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.971'
$ws.Range('E17').Value = '  +3.12%  '

$ws.Range('D18').Value = '51.526.27'
$ws.Range('E18').Value = '  +0.83%  '

$ws.Range('E19').Value = '  +2.88%  '

$ws.Range('E20').Value = '  +2.86%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.97'
$ws.Range('E21').Value = '  +1.45%  '

$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +2.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.50'
$ws.Range('E23').Value = '  +1.68%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.94'
$ws.Range('E24').Value = '  +0.73%  '

$ws.Range('E25').Value = '  +5.86%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.63'
$ws.Range('E26').Value = '  +25.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.80'
$ws.Range('E27').Value = '  +11.57%  '

$ws.Range('E28').Value = '  +0.22%  '

$ws.Range('E29').Value = '  +9.55%  '

$ws.Range('E30').Value = '  +0.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.95'
$ws.Range('E31').Value = '  +0.96%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.95'
$ws.Range('E32').Value = '  +0.56%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.43'
$ws.Range('E33').Value = '  +2.43%  '

$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.09'
$ws.Range('E34').Value = '  -1.96%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.13'
$ws.Range('E35').Value = '  +0.71%  '

$ws.Range('E36').Value = '  +6.05%  '

$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.27'
$ws.Range('E39').Value = '  +1.10%  '

$ws.Range('E40').Value = '  -2.14%  '

$ws.Range('E41').Value = '  +0.39%  '

$ws.Range('E42').Value = '  +2.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '125.07'
$ws.Range('E43').Value = '  +4.32%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.08'
$ws.Range('E44').Value = '  -0.35%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.290'
$ws.Range('E45').Value = '  +21.90%  '

$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('E47').Value = '  +2.96%  '

$ws.Range('D48').Value = '2.048.45'
$ws.Range('E48').Value = '  +1.35%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.25'
$ws.Range('E49').Value = '  +2.73%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0344'
$ws.Range('E50').Value = '  +10.45%  '

$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.15'
$ws.Range('E51').Value = '  +2.51%  '
